# Fruta / hortaliza, semanal
#
# A new weekly price record is inserted as row 148 ("Angeleno", 12-12-2022,
# Región Metropolitana). Every existing record from the old row 148 through
# the old row 249 shifts down by one row, and the record that used to be in
# row 249 now lands in the brand new row 250. The descriptive columns
# (A,B,C,E,F,G,H,I,J: Mercado/Región/Producto metadata) are identical for
# every row in this block, so only the "variable" columns
# (D,K,L,M,N,O,P,Q,R,S,T) actually need to move.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that hold the per-record data which shifts down one row at a time:
# D=Fecha, K=Variedad, L=Calidad, M=Volumen, N=Precio minimo, O=Precio maximo,
# P=Precio promedio ponderado, Q=Unidad de comercializacion, R=Origen,
# S=Precio $/Kg, T=Kg / unidad
$cols = @(4, 11, 12, 13, 14, 15, 16, 17, 18, 19, 20)

# Shift rows 148..249 down into rows 149..250, working from the bottom up so
# we never clobber a value before it has been copied.
for ($r = 250; $r -ge 149; $r--) {
    foreach ($c in $cols) {
        $srcVal = $ws.Cells.Item($r - 1, $c).Value2
        $ws.Cells.Item($r, $c).Value2 = $srcVal
    }
}

# Row 250 is brand new, so it needs the (constant-across-the-block)
# descriptive columns copied over explicitly.
$ws.Cells.Item(250, 1).Value2 = 4
$ws.Cells.Item(250, 2).Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(250, 3).Value2 = "Los Lagos"
$ws.Cells.Item(250, 5).Value2 = 10
$ws.Cells.Item(250, 6).Value2 = "Fruta"
$ws.Cells.Item(250, 7).Value2 = 100103
$ws.Cells.Item(250, 8).Value2 = "Frutos de hueso (carozo)"
$ws.Cells.Item(250, 9).Value2 = 100103002
$ws.Cells.Item(250, 10).Value2 = "Ciruela"

# Column D (Fecha) carries a date number-format throughout the sheet; since
# row 250 did not exist before, copy that format onto the new date cell too.
$ws.Cells.Item(250, 4).NumberFormat = $ws.Cells.Item(249, 4).NumberFormat

# Finally, write the brand new record into row 148.
$ws.Cells.Item(148, 4).Value2 = 44907
$ws.Cells.Item(148, 11).Value2 = "Angeleno"
$ws.Cells.Item(148, 12).Value2 = "Primera"
$ws.Cells.Item(148, 13).Value2 = 400
$ws.Cells.Item(148, 14).Value2 = 17000
$ws.Cells.Item(148, 15).Value2 = 18000
$ws.Cells.Item(148, 16).Value2 = 17500
$ws.Cells.Item(148, 17).Value2 = "$/caja 14 kilos granel"
$ws.Cells.Item(148, 18).Value2 = "Región Metropolitana"
$ws.Cells.Item(148, 19).Value2 = 1250
$ws.Cells.Item(148, 20).Value2 = 14

$dim = $ws.UsedRange.Address()
Write-Host "Final used range:" $dim
